# Black-Scholes Sheet1: add a new "T" (third) scenario in column F,
# mirroring the inputs/outputs already present in columns C (S) and D (T).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Inputs (F2:F7) ------------------------------------------------------
$ws1.Range("F2").Value = 10
$ws1.Range("F3").Value = 10
$ws1.Range("F4").Value = 5

$ws1.Range("F5").Value = 0.45
$ws1.Range("F5").NumberFormat = $ws1.Range("E5").NumberFormat()

$ws1.Range("F6").Value = 0.03
$ws1.Range("F6").NumberFormat = $ws1.Range("E6").NumberFormat()

$ws1.Range("F7").Value = 0.01
$ws1.Range("F7").NumberFormat = $ws1.Range("E7").NumberFormat()

# --- Formulas (F9:F18), matching the pattern used in columns C/D/E ------
$ws1.Range("F9").Formula = "=(LN(F2/F3)+(F6-F7+1/2*F5^2)*F4)/(F5*SQRT(F4))"
$ws1.Range("F9").NumberFormat = $ws1.Range("E9").NumberFormat()

$ws1.Range("F10").Formula = "=F9-(F5*SQRT(F4))"
$ws1.Range("F10").NumberFormat = $ws1.Range("E10").NumberFormat()

$ws1.Range("F11").NumberFormat = $ws1.Range("E11").NumberFormat()

$ws1.Range("F12").Formula = "=NORMSDIST(F9)"
$ws1.Range("F12").NumberFormat = $ws1.Range("E12").NumberFormat()

$ws1.Range("F13").Formula = "=NORMSDIST(F10)"
$ws1.Range("F13").NumberFormat = $ws1.Range("E13").NumberFormat()

$ws1.Range("F14").Formula = "=F2*EXP(-F7*F4)*F12-F3*EXP(-F6*F4)*F13"
$ws1.Range("F14").NumberFormat = $ws1.Range("E14").NumberFormat()

$ws1.Range("F15").NumberFormat = $ws1.Range("E15").NumberFormat()

$ws1.Range("F16").Formula = "=NORMSDIST(-F9)"
$ws1.Range("F16").NumberFormat = $ws1.Range("E16").NumberFormat()

$ws1.Range("F17").Formula = "=NORMSDIST(-F10)"
$ws1.Range("F17").NumberFormat = $ws1.Range("E17").NumberFormat()

$ws1.Range("F18").Formula = "=F3*EXP(-F6*F4)*F17-F2*EXP(-F7*F4)*F16"
$ws1.Range("F18").NumberFormat = $ws1.Range("E18").NumberFormat()

# --- Column F width, matching the C:E block ------------------------------
$ws1.Columns.Item(6).ColumnWidth = $ws1.Columns.Item(5).ColumnWidth()

# --- Selection / active sheet: Sheet1 becomes the active tab, cell F8 ---
$ws1.Range("F8").Select()
